$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.176.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "'1.902.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'306.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.5246"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("D8").Value = "'0.3770"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").Value = "'0.07268"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").Value = "'21.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").Value = "'0.8978"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").Value = "'0.08403"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.05%  "
$ws.Range("D13").Value = "'1.913.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").Value = "'94.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "'0.000008608"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "'14.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'27.209.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'5.058"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'2.139.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'10.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "'6.431"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").Value = "'146.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").Value = "'2.277"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.85%  "
$ws.Range("D27").Value = "'1.756"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "'114.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'4.931"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'4.792"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'0.09290"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "'0.8106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.31%  "
$ws.Range("D34").Value = "'0.05059"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'1.234"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.14%  "
$ws.Range("D36").Value = "'2.952"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.41%  "
$ws.Range("D37").Value = "'3.345"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("D38").Value = "'2.607"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("D39").Value = "'0.5708"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("D40").Value = "'0.01989"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").Value = "'6.665"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").Value = "'8.950"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("D44").Value = "'118.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "'0.1513"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "'0.4835"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").Value = "'1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "'1.614"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("D50").Value = "'37.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "'63.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
